$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3765.8333
$ws.Range("I74").Value = 797.5
$ws.Range("K74").Value = 797.5
$ws.Range("M74").Value = 138.5
$ws.Range("H77").Value = 3765.8333
$ws.Range("I77").Value = 797.5
$ws.Range("K77").Value = 3987.5
$ws.Range("M77").Value = 692.5
$ws.Range("H137").Value = 54410.26
$ws.Range("I137").Value = 1903.2858
$ws.Range("K137").Value = 5709.857400000001
$ws.Range("M137").Value = -3159.857400000001

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2406.2222
$ws.Range("I2").Value = 2214.9333
$ws.Range("J2").Value = 3362.6667
$ws.Range("K2").Value = 2214.9333
$ws.Range("L2").Value = 3362.6667
$ws.Range("M2").Value = -2101.9333
$ws.Range("N2").Value = -3588.6667
$ws.Range("H45").Value = 20379.117
$ws.Range("I45").Value = 20403.727
$ws.Range("J45").Value = 20334
$ws.Range("K45").Value = 20403.727
$ws.Range("L45").Value = 20334
$ws.Range("M45").Value = -20026.727
$ws.Range("N45").Value = -21088
$ws.Range("H116").Value = 2406.2222
$ws.Range("I116").Value = 2214.9333
$ws.Range("J116").Value = 3362.6667
$ws.Range("K116").Value = 2214.9333
$ws.Range("L116").Value = 3362.6667
$ws.Range("M116").Value = 79.06669999999986
$ws.Range("N116").Value = -7950.6667

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2406.2222
$ws.Range("I3").Value = 2214.9333
$ws.Range("J3").Value = 3362.6667
$ws.Range("K3").Value = 2214.9333
$ws.Range("L3").Value = 3362.6667
$ws.Range("M3").Value = -2100.9333
$ws.Range("N3").Value = -3590.6667
$ws.Range("H54").Value = 3900
$ws.Range("I54").Value = 3900
$ws.Range("K54").Value = 3900
$ws.Range("M54").Value = -3416
$ws.Range("H134").Value = 1753
$ws.Range("I134").Value = 1617.5
$ws.Range("K134").Value = 4852.5
$ws.Range("M134").Value = -2317.5

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 485.90475
$ws.Range("I7").Value = 308.875
$ws.Range("J7").Value = 594.8461
$ws.Range("K7").Value = 308.875
$ws.Range("L7").Value = 594.8461
$ws.Range("M7").Value = -195.875
$ws.Range("N7").Value = -820.8461
$ws.Range("H22").Value = 462.5
$ws.Range("I22").Value = 425
$ws.Range("K22").Value = 425
$ws.Range("M22").Value = -75
$ws.Range("H51").Value = 11999.3
$ws.Range("J51").Value = 11999.3
$ws.Range("L51").Value = 11999.3
$ws.Range("N51").Value = -13471.3
$ws.Range("H59").Value = 16108.889
$ws.Range("J59").Value = 16108.889
$ws.Range("L59").Value = 16108.889
$ws.Range("N59").Value = -18398.889
$ws.Range("H60").Value = 9727.571
$ws.Range("J60").Value = 10000
$ws.Range("L60").Value = 10000
$ws.Range("N60").Value = -11022
$ws.Range("H61").Value = 11999.3
$ws.Range("J61").Value = 11999.3
$ws.Range("L61").Value = 11999.3
$ws.Range("N61").Value = -12695.3
$ws.Range("H99").Value = 7150
$ws.Range("I99").Value = 800
$ws.Range("K99").Value = 800
$ws.Range("M99").Value = 698
$ws.Range("H107").Value = 972.6316
$ws.Range("I107").Value = 668.9167
$ws.Range("J107").Value = 1493.2858
$ws.Range("K107").Value = 668.9167
$ws.Range("L107").Value = 1493.2858
$ws.Range("M107").Value = 1251.0833
$ws.Range("N107").Value = -5333.2858
$ws.Range("H126").Value = 7150
$ws.Range("I126").Value = 800
$ws.Range("K126").Value = 2400
$ws.Range("M126").Value = 70
$ws.Range("H134").Value = 2378.8
$ws.Range("J134").Value = 3405.2
$ws.Range("L134").Value = 10215.6
$ws.Range("N134").Value = -15285.6

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 322.47827
$ws.Range("I23").Value = 331.33334
$ws.Range("J23").Value = 319.35294
$ws.Range("K23").Value = 994.0000200000001
$ws.Range("L23").Value = 958.05882
$ws.Range("M23").Value = -759.0000200000001
$ws.Range("N23").Value = -1428.05882
$ws.Range("H107").Value = 1385.7142
$ws.Range("J107").Value = 1366.6666
$ws.Range("L107").Value = 4099.9998
$ws.Range("N107").Value = -7939.9998

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 306.47058
$ws.Range("J2").Value = 364.4
$ws.Range("L2").Value = 364.4
$ws.Range("N2").Value = -590.4
$ws.Range("H80").Value = 8215.956
$ws.Range("I80").Value = 9167.5
$ws.Range("J80").Value = 6041
$ws.Range("K80").Value = 9167.5
$ws.Range("L80").Value = 6041
$ws.Range("M80").Value = -8169.5
$ws.Range("N80").Value = -8037
$ws.Range("H83").Value = 8215.956
$ws.Range("I83").Value = 9167.5
$ws.Range("J83").Value = 6041
$ws.Range("K83").Value = 45837.5
$ws.Range("L83").Value = 30205
$ws.Range("M83").Value = -40845.5
$ws.Range("N83").Value = -40189

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H47").Value = 11065
$ws.Range("J47").Value = 11065
$ws.Range("L47").Value = 11065
$ws.Range("N47").Value = -12045
$ws.Range("H50").Value = 221247.5
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 221247.5
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 221247.5
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -222521.5
$ws.Range("H52").Value = 11065
$ws.Range("J52").Value = 11065
$ws.Range("L52").Value = 11065
$ws.Range("N52").Value = -11531
$ws.Range("H54").Value = 31247.5
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 31247.5
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 31247.5
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -32535.5
$ws.Range("H56").Value = 19000
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 19000
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 19000
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -20382
$ws.Range("H58").Value = 12650.375
$ws.Range("J58").Value = 20102.5
$ws.Range("L58").Value = 20102.5
$ws.Range("N58").Value = -20622.5
$ws.Range("H93").Value = 38445.555
$ws.Range("I93").Value = 1516.3334
$ws.Range("K93").Value = 1516.3334
$ws.Range("M93").Value = -268.3334
$ws.Range("H136").Value = 4036.4285
$ws.Range("I136").Value = 2167
$ws.Range("K136").Value = 6501
$ws.Range("M136").Value = -3951

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 20949.75
$ws.Range("I4").Value = 2399
$ws.Range("K4").Value = 2399
$ws.Range("M4").Value = -2286
$ws.Range("H81").Value = 2239
$ws.Range("I81").Value = 2239
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 4478
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -3417
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 2239
$ws.Range("I84").Value = 2239
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 22390
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -17086
$ws.Range("N84").ClearContents()
